# HAJ-212 / HAJJ-206: align sheet name for the applicant ritual data template
# with the applicant health data template ("Applicant Health Data" -> "Applicant Ritual Data").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Applicant Ritual Data"
